$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 92
$ws.Range("H92").Value = 618.6667
$ws.Range("I92").Value = 677.64703
$ws.Range("J92").Value = 475.42856
$ws.Range("K92").Value = 677.64703
$ws.Range("L92").Value = 475.42856
$ws.Range("M92").Value = 570.35297
$ws.Range("N92").Value = -2971.42856
# Row 94
$ws.Range("H94").Value = 1871.4286
$ws.Range("J94").Value = 1933.3334
$ws.Range("L94").Value = 1933.3334
$ws.Range("N94").Value = -2835.3334
# Row 100
$ws.Range("H100").Value = 2391.05
$ws.Range("I100").Value = 2573.389
$ws.Range("J100").Value = 750
$ws.Range("K100").Value = 2573.389
$ws.Range("L100").Value = 750
$ws.Range("M100").Value = -2032.389
$ws.Range("N100").Value = -1832
# Row 103
$ws.Range("H103").Value = 0
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 0
$ws.Range("L103").ClearContents()
$ws.Range("M103").ClearContents()
$ws.Range("N103").Value = 0
# Row 107
$ws.Range("H107").Value = 1223.7858
$ws.Range("I107").Value = 1399.8334
$ws.Range("J107").Value = 167.5
$ws.Range("K107").Value = 1399.8334
$ws.Range("L107").Value = 167.5
$ws.Range("M107").Value = 520.1666
$ws.Range("N107").Value = -4007.5
# Row 131
$ws.Range("H131").Value = 955.8
$ws.Range("I131").Value = 930
$ws.Range("J131").Value = 994.5
$ws.Range("K131").Value = 2790
$ws.Range("L131").Value = 2983.5
$ws.Range("M131").Value = 2250
$ws.Range("N131").Value = -13063.5
# Row 138
$ws.Range("H138").Value = 3392.6858
$ws.Range("J138").Value = 3851.7778
$ws.Range("L138").Value = 11555.3334
$ws.Range("N138").Value = -21835.3334

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 102
$ws.Range("H102").Value = 968.8461
$ws.Range("I102").Value = 968.8461
$ws.Range("K102").Value = 968.8461
$ws.Range("M102").Value = 653.1539
# Row 110
$ws.Range("H110").Value = 2626.4
$ws.Range("I110").Value = 1028.2858
$ws.Range("J110").Value = 25000
$ws.Range("K110").Value = 1028.2858
$ws.Range("L110").Value = 25000
$ws.Range("M110").Value = 1016.7142
$ws.Range("N110").Value = -29090
# Row 132
$ws.Range("H132").Value = 1917.2222
$ws.Range("I132").Value = 1925.6923
$ws.Range("K132").Value = 5777.0769
$ws.Range("M132").Value = -3247.0769

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 2757.739
$ws.Range("I86").Value = 1578.1177
$ws.Range("K86").Value = 1578.1177
$ws.Range("M86").Value = -455.1177
# Row 89
$ws.Range("H89").Value = 2757.739
$ws.Range("I89").Value = 1578.1177
$ws.Range("K89").Value = 7890.5885
$ws.Range("M89").Value = -2274.5885
# Row 94
$ws.Range("H94").Value = 3331.7
$ws.Range("I94").Value = 3978.75
$ws.Range("K94").Value = 3978.75
$ws.Range("M94").Value = -3527.75
# Row 105
$ws.Range("H105").Value = 2682.7144
$ws.Range("I105").Value = 2713.5
$ws.Range("K105").Value = 2713.5
$ws.Range("M105").Value = -966.5

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 189.57143
$ws.Range("I7").Value = 232.5
$ws.Range("J7").Value = 132.33333
$ws.Range("K7").Value = 232.5
$ws.Range("L7").Value = 132.33333
$ws.Range("M7").Value = -119.5
$ws.Range("N7").Value = -358.33333
# Row 16
$ws.Range("H16").Value = 5011
$ws.Range("I16").Value = 5011
$ws.Range("K16").Value = 5011
$ws.Range("M16").Value = -4724
# Row 31
$ws.Range("H31").Value = 2192.8
$ws.Range("I31").Value = 1994.6666
$ws.Range("K31").Value = 1994.6666
$ws.Range("M31").Value = -1699.6666
# Row 34
$ws.Range("H34").Value = 2192.8
$ws.Range("I34").Value = 1994.6666
$ws.Range("K34").Value = 1994.6666
$ws.Range("M34").Value = -1792.6666
# Row 62
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").ClearContents()
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = 0
# Row 65
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").ClearContents()
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = 0
# Row 113
$ws.Range("H113").Value = 5011
$ws.Range("I113").Value = 5011
$ws.Range("K113").Value = 5011
$ws.Range("M113").Value = -2841

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 7
$ws.Range("H7").Value = 150
$ws.Range("I7").Value = 150
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 450
$ws.Range("L7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = -338
# Row 80
$ws.Range("H80").Value = 11090.583
$ws.Range("I80").Value = 4771.75
$ws.Range("J80").Value = 14250
$ws.Range("K80").Value = 14315.25
$ws.Range("L80").Value = 42750
$ws.Range("M80").Value = -13379.25
$ws.Range("N80").Value = -44622
# Row 83
$ws.Range("H83").Value = 11090.583
$ws.Range("I83").Value = 4771.75
$ws.Range("J83").Value = 14250
$ws.Range("K83").Value = 42945.75
$ws.Range("L83").Value = 128250
$ws.Range("M83").Value = -38265.75
$ws.Range("N83").Value = -137610
# Row 92
$ws.Range("H92").Value = 5191
$ws.Range("I92").Value = 379
$ws.Range("K92").Value = 1137
$ws.Range("M92").Value = 111

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 2939.5557
$ws.Range("I80").Value = 2598.8
$ws.Range("J80").Value = 3365.5
$ws.Range("K80").Value = 2598.8
$ws.Range("L80").Value = 3365.5
$ws.Range("M80").Value = -1600.8
$ws.Range("N80").Value = -5361.5
# Row 83
$ws.Range("H83").Value = 2939.5557
$ws.Range("I83").Value = 2598.8
$ws.Range("J83").Value = 3365.5
$ws.Range("K83").Value = 12994
$ws.Range("L83").Value = 16827.5
$ws.Range("M83").Value = -8002
$ws.Range("N83").Value = -26811.5
# Row 97
$ws.Range("H97").Value = 1819.8
$ws.Range("I97").Value = 2037.25
$ws.Range("J97").Value = 950
$ws.Range("K97").Value = 2037.25
$ws.Range("L97").Value = 950
$ws.Range("M97").Value = -1541.25
$ws.Range("N97").Value = -1942

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 82
$ws.Range("H82").Value = 2049.2856
$ws.Range("J82").Value = 2976.5
$ws.Range("L82").Value = 2976.5
$ws.Range("N82").Value = -3698.5
# Row 85
$ws.Range("H85").Value = 2049.2856
$ws.Range("J85").Value = 2976.5
$ws.Range("L85").Value = 2976.5
$ws.Range("N85").Value = -5472.5
# Row 99
$ws.Range("H99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("L99").ClearContents()
$ws.Range("N99").Value = 0
# Row 132
$ws.Range("H132").Value = 5832.25
$ws.Range("J132").Value = 5832.25
$ws.Range("L132").Value = 17496.75
$ws.Range("N132").Value = -22556.75

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 4881.7144
$ws.Range("I62").Value = 3533
$ws.Range("K62").Value = 3533
$ws.Range("M62").Value = -2909
# Row 65
$ws.Range("H65").Value = 4881.7144
$ws.Range("I65").Value = 3533
$ws.Range("K65").Value = 17665
$ws.Range("M65").Value = -14545
# Row 96
$ws.Range("H96").Value = 1899.5
$ws.Range("I96").Value = 1899.5
$ws.Range("K96").Value = 1899.5
$ws.Range("M96").Value = -526.5

Write-Host "Done applying changes"